# Automatische test-sync: 2025-07-23 18:39:50
# Append a new mail-log row (row 7) to the "Logs" sheet, extend the
# conditional-formatting ranges to cover it, and bump the matching tally
# on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 7

$logs.Cells.Item($newRow, 1).Value = "Wat zijn jullie openingstijden?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #1: Wat zijn jullie openingstijden?"
$logs.Cells.Item($newRow, 4).Value = "Openingstijden / Locatie"
$logs.Cells.Item($newRow, 5).Value = "Beste klant,`nDank u wel voor uw interesse in onze diensten. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Cells.Item($newRow, 6).Value = "2025-07-23 18:39:49"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Nee"
$logs.Cells.Item($newRow, 9).Value = "Ja"
$logs.Cells.Item($newRow, 10).Value = "Ja"

# Writing multi-line text into a brand-new row auto-expands its height;
# put it back to the sheet's standard height (matches the other rows,
# none of which carry an explicit row height).
$logs.Rows.Item($newRow).AutoFit()

# The conditional-formatting bands on columns D, G, H, I and J used to stop
# at row 6 (the old last data row) - stretch each one down to the new last
# row (7) while keeping their rules/order/priority untouched.
function Extend-FormatConditions($columnLetter) {
    $oldRange = $logs.Range($columnLetter + "2:" + $columnLetter + "6")
    $newRange = $logs.Range($columnLetter + "2:" + $columnLetter + "7")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

Extend-FormatConditions "D"
Extend-FormatConditions "G"
Extend-FormatConditions "H"
Extend-FormatConditions "I"
Extend-FormatConditions "J"

# Dashboard tally: "Openingstijden / Locatie" count goes from 5 to 6.
$dashboard.Cells.Item(2, 2).Value = 6
